# 1) Merge the three runs "TCP(" + "dport" + "=[80,443])" (with the two
#    spell-check proofErr markers between them) into a single run
#    "TCP(dport=[80,443])", keeping the formatting of the first run/the
#    paragraph mark (Word naturally collapses proofErr + merges runs when
#    the found text is replaced in one shot).
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("TCP(dport=[80,443])", $false, $false, $false, $false, $false, $true, 1, $false, "TCP(dport=[80,443])", 2)
if (-not $found) {
    throw "Could not find the TCP(dport=[80,443]) text to normalize into a single run"
}

# 2) Append two new paragraphs at the very end of the document (after the
#    sendp(p, iface="veth11") line): a "=====" divider line and a new
#    "sudo tcpdump -n -i veth30" usage example, matching the look & feel of
#    the other divider/command pairs already used throughout the doc.
$endRange = $d.Content
$endRange.Collapse(0)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:left w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:bottom w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:right w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:wordWrap w:val="0"/><w:spacing w:after="150" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/></w:rPr><w:t>=====================================================================================</w:t></w:r></w:p><w:p><w:pPr><w:pBdr><w:top w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:left w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:bottom w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/><w:right w:val="single" w:sz="6" w:space="7" w:color="CCCCCC"/></w:pBdr><w:shd w:val="clear" w:color="auto" w:fill="F5F5F5"/><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:wordWrap w:val="0"/><w:spacing w:after="150" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="vi-VN"/></w:rPr><w:t>sudo tcpdump -n -i veth3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:eastAsia="Times New Roman" w:hAnsi="Consolas" w:cs="Courier New"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/></w:rPr><w:t>0</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($xml)

